$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7985977530479431
$ws.Range("B1").Value = 2.468005895614624
$ws.Range("C1").Value = 1.424812436103821
$ws.Range("D1").Value = 1.045518159866333
$ws.Range("E1").Value = 1.212092280387878
